$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties in AC1:AE1
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold / border / centered) from an existing header cell (AB1) to the new headers
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the team record values (Wins=51, Losses=64, Ties=0) for every data row (2-44)
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 51   # AC
    $ws.Cells.Item($r, 30).Value = 64   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}
